$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.434937333333333
$ws.Range("H2").Value = 4.304812
$ws.Range("I2").Value = 0.5010808920723563
$ws.Range("J2").Value = 0.5010808920723562
$ws.Range("M2").Value = 15.50220733333333
$ws.Range("N2").Value = 46.506622
$ws.Range("O2").Value = 0.5994675913188158
$ws.Range("P2").Value = 0.5994675913188158
$ws.Range("Q2").Value = 22.24469605167378
$ws.Range("R2").Value = 200.202264465064
$ws.Range("S2").Value = 0.3003817554264989
$ws.Range("T2").Value = 0.3003817554264989
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.434937333333333
$ws.Range("H3").Value = 4.304812
$ws.Range("I3").Value = 0.5010808920723563
$ws.Range("J3").Value = 0.5010808920723562
$ws.Range("O3").Value = 0.04399860030713892
$ws.Range("P3").Value = 0.04399860030713892
$ws.Range("Q3").Value = 1.632674567741333
$ws.Range("R3").Value = 14.694071109672
$ws.Range("S3").Value = 0.02204685789183622
$ws.Range("T3").Value = 0.02204685789183622
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.434937333333333
$ws.Range("H4").Value = 4.304812
$ws.Range("I4").Value = 0.5010808920723563
$ws.Range("J4").Value = 0.5010808920723562
$ws.Range("M4").Value = 8.848210666666667
$ws.Range("N4").Value = 26.544632
$ws.Range("O4").Value = 0.3421587275782868
$ws.Range("P4").Value = 0.3421587275782868
$ws.Range("Q4").Value = 12.69662781879822
$ws.Range("R4").Value = 114.269650369184
$ws.Range("S4").Value = 0.1714492004452703
$ws.Range("T4").Value = 0.1714492004452703
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.434937333333333
$ws.Range("H5").Value = 4.304812
$ws.Range("I5").Value = 0.5010808920723563
$ws.Range("J5").Value = 0.5010808920723562
$ws.Range("M5").Value = 0.371739
$ws.Range("N5").Value = 1.115217
$ws.Range("O5").Value = 0.01437508079575842
$ws.Range("P5").Value = 0.01437508079575841
$ws.Range("Q5").Value = 0.533422169356
$ws.Range("R5").Value = 4.800799524204
$ws.Range("S5").Value = 0.007203078308750824
$ws.Range("T5").Value = 0.007203078308750821
$ws.Range("G6").Value = 0.9964423333333334
$ws.Range("I6").Value = 0.3479582011609289
$ws.Range("J6").Value = 0.3479582011609288
$ws.Range("M6").Value = 15.50220733333333
$ws.Range("N6").Value = 46.506622
$ws.Range("O6").Value = 0.5994675913188158
$ws.Range("P6").Value = 0.5994675913188158
$ws.Range("Q6").Value = 15.44705564704378
$ws.Range("R6").Value = 139.023500823394
$ws.Range("S6").Value = 0.20858966472957
$ws.Range("T6").Value = 0.2085896647295699
$ws.Range("G7").Value = 0.9964423333333334
$ws.Range("I7").Value = 0.3479582011609289
$ws.Range("J7").Value = 0.3479582011609288
$ws.Range("O7").Value = 0.04399860030713892
$ws.Range("P7").Value = 0.04399860030713892
$ws.Range("S7").Value = 0.01530967381647075
$ws.Range("T7").Value = 0.01530967381647075
$ws.Range("G8").Value = 0.9964423333333334
$ws.Range("I8").Value = 0.3479582011609289
$ws.Range("J8").Value = 0.3479582011609288
$ws.Range("M8").Value = 8.848210666666667
$ws.Range("N8").Value = 26.544632
$ws.Range("O8").Value = 0.3421587275782868
$ws.Range("P8").Value = 0.3421587275782868
$ws.Range("Q8").Value = 8.816731682518222
$ws.Range("R8").Value = 79.350585142664
$ws.Range("S8").Value = 0.119056935359653
$ws.Range("T8").Value = 0.1190569353596529
$ws.Range("G9").Value = 0.9964423333333334
$ws.Range("I9").Value = 0.3479582011609289
$ws.Range("J9").Value = 0.3479582011609288
$ws.Range("M9").Value = 0.371739
$ws.Range("N9").Value = 1.115217
$ws.Range("O9").Value = 0.01437508079575842
$ws.Range("P9").Value = 0.01437508079575841
$ws.Range("Q9").Value = 0.370416476551
$ws.Range("R9").Value = 3.333748288959
$ws.Range("S9").Value = 0.005001927255235113
$ws.Range("T9").Value = 0.00500192725523511
$ws.Range("G10").Value = 0.4323043333333333
$ws.Range("H10").Value = 1.296913
$ws.Range("I10").Value = 0.150960906766715
$ws.Range("J10").Value = 0.1509609067667149
$ws.Range("M10").Value = 15.50220733333333
$ws.Range("N10").Value = 46.506622
$ws.Range("O10").Value = 0.5994675913188158
$ws.Range("P10").Value = 0.5994675913188158
$ws.Range("Q10").Value = 6.701671406431778
$ws.Range("R10").Value = 60.315042657886
$ws.Range("S10").Value = 0.09049617116274694
$ws.Range("T10").Value = 0.09049617116274691
$ws.Range("G11").Value = 0.4323043333333333
$ws.Range("H11").Value = 1.296913
$ws.Range("I11").Value = 0.150960906766715
$ws.Range("J11").Value = 0.1509609067667149
$ws.Range("O11").Value = 0.04399860030713892
$ws.Range("P11").Value = 0.04399860030713892
$ws.Range("Q11").Value = 0.4918767350753334
$ws.Range("R11").Value = 4.426890615678
$ws.Range("S11").Value = 0.006642068598831956
$ws.Range("T11").Value = 0.006642068598831954
$ws.Range("G12").Value = 0.4323043333333333
$ws.Range("H12").Value = 1.296913
$ws.Range("I12").Value = 0.150960906766715
$ws.Range("J12").Value = 0.1509609067667149
$ws.Range("M12").Value = 8.848210666666667
$ws.Range("N12").Value = 26.544632
$ws.Range("O12").Value = 0.3421587275782868
$ws.Range("P12").Value = 0.3421587275782868
$ws.Range("Q12").Value = 3.825119813446222
$ws.Range("R12").Value = 34.426078321016
$ws.Range("S12").Value = 0.05165259177336358
$ws.Range("T12").Value = 0.05165259177336357
$ws.Range("G13").Value = 0.4323043333333333
$ws.Range("H13").Value = 1.296913
$ws.Range("I13").Value = 0.150960906766715
$ws.Range("J13").Value = 0.1509609067667149
$ws.Range("M13").Value = 0.371739
$ws.Range("N13").Value = 1.115217
$ws.Range("O13").Value = 0.01437508079575842
$ws.Range("P13").Value = 0.01437508079575841
$ws.Range("Q13").Value = 0.160704380569
$ws.Range("R13").Value = 1.446339425121
$ws.Range("S13").Value = 0.002170075231772481
$ws.Range("T13").Value = 0.00217007523177248
